$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.056.20"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "2.240.37"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'246.78"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").Value = "'75.82"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "'40.46"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  -3.88%  "
$ws.Range("D12").Value = "'7.18"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "2.575.88"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "'14.87"
$ws.Range("E15").Value = "  -4.15%  "
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "2.271.64"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "41.926.84"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "'6.14"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "'71.66"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "'2.26"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'231.11"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -5.87%  "
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  +11.36%  "
$ws.Range("D29").Value = "'2.15"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "'168.87"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "'20.55"
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("D32").Value = "'33.47"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").Value = "'0.126"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "'4.51"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "'4.85"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "'13.41"
$ws.Range("E39").Value = "  -7.95%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").Value = "'2.18"
$ws.Range("E41").Value = "  -7.35%  "
$ws.Range("D42").Value = "'111.50"
$ws.Range("E42").Value = "  +13.14%  "
$ws.Range("D43").Value = "'0.203"
$ws.Range("E43").Value = "  -5.37%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").Value = "'8.74"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").Value = "'0.996"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("E49").Value = "  -11.23%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'4.18"
$ws.Range("E51").Value = "  -2.72%  "
